# Fruta / hortaliza, semanal
# Insert a new weekly record at row 454 of the "Ajo" sheet, pushing the
# existing rows 454:487 down to 455:488 (dimension grows from R487 to R488).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 454 (shifts 454:487 -> 455:488).
$ws.Rows.Item(454).Insert()

# Populate the newly inserted row 454 with the new weekly price record.
$ws.Cells.Item(454, 1).Value = 8
$ws.Cells.Item(454, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(454, 3).Value = "Coquimbo"
$ws.Cells.Item(454, 4).Value = 45106
$ws.Cells.Item(454, 5).Value = 4
$ws.Cells.Item(454, 6).Value = 100112003
$ws.Cells.Item(454, 7).Value = "Ajo"
$ws.Cells.Item(454, 8).Value = "Chino"
$ws.Cells.Item(454, 9).Value = "Primera"
$ws.Cells.Item(454, 10).Value = 300
$ws.Cells.Item(454, 11).Value = 17500
$ws.Cells.Item(454, 12).Value = 18000
$ws.Cells.Item(454, 13).Value = 17750
$ws.Cells.Item(454, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(454, 15).Value = "China"
$ws.Cells.Item(454, 16).Value = 1775
$ws.Cells.Item(454, 17).Value = 10
$ws.Cells.Item(454, 18).Value = "Hortaliza"
